$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 402
$ws1.Range("F7").Value = 878
$ws1.Range("F8").Value = 59
$ws1.Range("F9").Value = 524
$ws1.Range("F11").Value = 298
$ws1.Range("F12").Value = 1157
$ws1.Range("F14").Value = 250
$ws1.Range("F15").Value = 38
$ws1.Range("F17").Value = 6682
$ws1.Range("F21").Value = 7603
$ws1.Range("F22").Value = 42
$ws1.Range("F24").Value = 3406
$ws1.Range("F26").Value = 2124
$ws1.Range("F27").Value = 901
$ws1.Range("F28").Value = 4519
$ws1.Range("F29").Value = 155
$ws1.Range("F31").Value = 72
$ws1.Range("F32").Value = 235
$ws1.Range("F34").Value = 1728
$ws1.Range("F36").Value = 181
$ws1.Range("F37").Value = 56
$ws1.Range("F39").Value = 1224
$ws1.Range("F40").Value = 1823
$ws1.Range("F41").Value = 2146

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 402
$ws4.Range("F9").Value = 878
$ws4.Range("F10").Value = 59
$ws4.Range("F11").Value = 524
$ws4.Range("F13").Value = 298
$ws4.Range("F14").Value = 1157
$ws4.Range("F17").Value = 250
$ws4.Range("F18").Value = 38
$ws4.Range("F20").Value = 6682
$ws4.Range("F24").Value = 7603
$ws4.Range("F25").Value = 42
$ws4.Range("F27").Value = 3406
$ws4.Range("F29").Value = 2124
$ws4.Range("F30").Value = 901
$ws4.Range("F31").Value = 4519
$ws4.Range("F32").Value = 155
$ws4.Range("F34").Value = 72
$ws4.Range("F36").Value = 235
$ws4.Range("F38").Value = 1728
$ws4.Range("F40").Value = 181
$ws4.Range("F41").Value = 56
$ws4.Range("F44").Value = 1224
$ws4.Range("F45").Value = 1823
$ws4.Range("F47").Value = 2146
